$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - AgregarArticulo
$ws.Range("B7").Value2 = "AgregarArticulo"
$ws.Range("C7").Value2 = "En el portal de myshopify agrego un articulo al carrito de compras"
$ws.Range("D7").Value2 = "Me encuentro en el home "
$ws.Range("E7").Value2 = "1. le doy click a la primera opción de los articulos 2. Cuando me manda a la información del articulo darle click al boton de agregar al carrito. 3. recargar la pagina 4. entrar al carrito."
$ws.Range("F7").Value2 = "El carrito tiene el articulo "

# Row 8 - EliminarArticulo
$ws.Range("B8").Value2 = "EliminarArticulo"
$ws.Range("C8").Value2 = "En el portal de myshopify elimino un articulo del carroto de compras "
$ws.Range("D8").Value2 = "Me encuentro en el home "
$ws.Range("E8").Value2 = "1. le doy click a la primera opción de los articulos 2. Cuando me manda a la información del articulo darle click al boton de agregar al carrito. 3. recargar la pagina 4. entrar al carrito. 5. Elimino el articulo "
$ws.Range("F8").Value2 = "El carrito esta vacio"

# Wrap text on C and E columns (matching the description/steps columns)
$ws.Range("C7").WrapText = $true
$ws.Range("E7").WrapText = $true
$ws.Range("C8").WrapText = $true
$ws.Range("E8").WrapText = $true

# Row heights
$ws.Rows.Item(7).RowHeight = 90
$ws.Rows.Item(8).RowHeight = 105

# Sheet view changes: remove frozen/scrolled topLeftCell (back to A1), change selection to G8
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("G8").Select()
